$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 51, shifting existing rows 51.. down by one.
$ws.Rows("51:51").Insert()

# Populate the newly inserted row 51 with the new record.
$ws.Range("A51").Value = 10
$ws.Range("B51").Value = "Vega Modelo de Temuco"
$ws.Range("C51").Value = "La Araucanía"
$ws.Range("D51").Value = 45259
$ws.Range("D51").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E51").Value = 9
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100101
$ws.Range("H51").Value = "Berries"
$ws.Range("I51").Value = 100101001
$ws.Range("J51").Value = "Arándano (blue)"
$ws.Range("K51").Value = "Sin especificar"
$ws.Range("L51").Value = "Primera"
$ws.Range("M51").Value = 50
$ws.Range("N51").Value = 4000
$ws.Range("O51").Value = 4000
$ws.Range("P51").Value = 4000
$ws.Range("Q51").Value = "$/kilo"
$ws.Range("R51").Value = "Región del Maule"
$ws.Range("S51").Value = 4000
$ws.Range("T51").Value = 1
